# edit.ps1 - Update cryptocurrency price table (Coinranking snapshot refresh)
# Commit: "Updated symbol list on Fri Feb  3 12:26:53 UTC 2023 with GitHub Actions"
#
# Columns D (Price), E (Volume 1h %), G (Hora/hour) are stored in the workbook
# as literal text (e.g. "327.95", "-0.58%", "11"), not as numbers/percentages.
# Forcing NumberFormat to "@" (Text) before assigning the new values keeps them
# stored the same way Excel would otherwise auto-convert numeric-looking text
# (e.g. "331.41") or percent-looking text (e.g. "0.41%") into real numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the ranges that receive new values as Text so assignment below
# preserves them as literal strings instead of being parsed as numbers.
$ws.Range("D2:D26").NumberFormat = "@"
$ws.Range("D39:D42").NumberFormat = "@"
$ws.Range("D44:D51").NumberFormat = "@"
$ws.Range("E2:E26").NumberFormat = "@"
$ws.Range("E39:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "331.41"
$ws.Range("E2").Value = "0.41%"
$ws.Range("G2").Value = "12"

# Row 3
$ws.Range("D3").Value = "39.68"
$ws.Range("E3").Value = "-1.20%"
$ws.Range("G3").Value = "12"

# Row 4
$ws.Range("D4").Value = "5.775"
$ws.Range("E4").Value = "5.30%"
$ws.Range("G4").Value = "12"

# Row 5
$ws.Range("D5").Value = "0.08059"
$ws.Range("E5").Value = "-0.80%"
$ws.Range("G5").Value = "12"

# Row 6
$ws.Range("D6").Value = "1.992"
$ws.Range("E6").Value = "3.52%"
$ws.Range("G6").Value = "12"

# Row 7
$ws.Range("D7").Value = "4.500"
$ws.Range("E7").Value = "-0.72%"
$ws.Range("G7").Value = "12"

# Row 8
$ws.Range("D8").Value = "8.643"
$ws.Range("E8").Value = "-0.32%"
$ws.Range("G8").Value = "12"

# Row 9
$ws.Range("D9").Value = "2.999"
$ws.Range("E9").Value = "0.73%"
$ws.Range("G9").Value = "12"

# Row 10
$ws.Range("D10").Value = "0.9219"
$ws.Range("E10").Value = "-2.87%"
$ws.Range("G10").Value = "12"

# Row 11
$ws.Range("D11").Value = "0.1271"
$ws.Range("E11").Value = "-3.57%"
$ws.Range("G11").Value = "12"

# Row 12
$ws.Range("D12").Value = "0.1955"
$ws.Range("E12").Value = "-1.93%"
$ws.Range("G12").Value = "12"

# Row 13
$ws.Range("D13").Value = "8.742"
$ws.Range("E13").Value = "20.15%"
$ws.Range("G13").Value = "12"

# Row 14
$ws.Range("D14").Value = "0.09240"
$ws.Range("E14").Value = "0.60%"
$ws.Range("G14").Value = "12"

# Row 15
$ws.Range("D15").Value = "0.03570"
$ws.Range("E15").Value = "0.09%"
$ws.Range("G15").Value = "12"

# Row 16
$ws.Range("D16").Value = "0.1052"
$ws.Range("E16").Value = "9.67%"
$ws.Range("G16").Value = "12"

# Row 17
$ws.Range("D17").Value = "0.001305"
$ws.Range("E17").Value = "-1.76%"
$ws.Range("G17").Value = "12"

# Row 18
$ws.Range("D18").Value = "0.006213"
$ws.Range("E18").Value = "-0.26%"
$ws.Range("G18").Value = "12"

# Row 19
$ws.Range("D19").Value = "3.368"
$ws.Range("E19").Value = "-0.04%"
$ws.Range("G19").Value = "12"

# Row 20
$ws.Range("D20").Value = "0.3480"
$ws.Range("E20").Value = "-0.95%"
$ws.Range("G20").Value = "12"

# Row 21
$ws.Range("D21").Value = "0.1361"
$ws.Range("E21").Value = "2.14%"
$ws.Range("G21").Value = "12"

# Row 22
$ws.Range("D22").Value = "0.2760"
$ws.Range("E22").Value = "12.64%"
$ws.Range("G22").Value = "12"

# Row 23
$ws.Range("D23").Value = "0.04399"
$ws.Range("E23").Value = "-0.57%"
$ws.Range("G23").Value = "12"

# Row 24
$ws.Range("D24").Value = "0.001260"
$ws.Range("E24").Value = "2.82%"
$ws.Range("G24").Value = "12"

# Row 25
$ws.Range("D25").Value = "0.004616"
$ws.Range("E25").Value = "6.87%"
$ws.Range("G25").Value = "12"

# Row 26
$ws.Range("D26").Value = "0.0001189"
$ws.Range("E26").Value = "-0.93%"
$ws.Range("G26").Value = "12"

# Row 27
$ws.Range("G27").Value = "12"

# Row 28
$ws.Range("G28").Value = "12"

# Row 29
$ws.Range("G29").Value = "12"

# Row 30
$ws.Range("G30").Value = "12"

# Row 31
$ws.Range("G31").Value = "12"

# Row 32
$ws.Range("G32").Value = "12"

# Row 33
$ws.Range("G33").Value = "12"

# Row 34
$ws.Range("G34").Value = "12"

# Row 35
$ws.Range("G35").Value = "12"

# Row 36
$ws.Range("G36").Value = "12"

# Row 37
$ws.Range("G37").Value = "12"

# Row 38
$ws.Range("G38").Value = "12"

# Row 39
$ws.Range("D39").Value = "0.02480"
$ws.Range("E39").Value = "-0.91%"
$ws.Range("G39").Value = "12"

# Row 40
$ws.Range("D40").Value = "0.05547"
$ws.Range("E40").Value = "5.37%"
$ws.Range("G40").Value = "12"

# Row 41
$ws.Range("D41").Value = "0.007446"
$ws.Range("E41").Value = "-4.16%"
$ws.Range("G41").Value = "12"

# Row 42
$ws.Range("D42").Value = "0.009942"
$ws.Range("E42").Value = "7.18%"
$ws.Range("G42").Value = "12"

# Row 43
$ws.Range("E43").Value = "-1.60%"
$ws.Range("G43").Value = "12"

# Row 44
$ws.Range("D44").Value = "0.002107"
$ws.Range("E44").Value = "-2.51%"
$ws.Range("G44").Value = "12"

# Row 45
$ws.Range("D45").Value = "0.01147"
$ws.Range("E45").Value = "18.93%"
$ws.Range("G45").Value = "12"

# Row 46
$ws.Range("D46").Value = "0.00006691"
$ws.Range("E46").Value = "1.50%"
$ws.Range("G46").Value = "12"

# Row 47
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").Value = "-0.06%"
$ws.Range("G47").Value = "12"

# Row 48
$ws.Range("D48").Value = "0.003030"
$ws.Range("E48").Value = "-6.33%"
$ws.Range("G48").Value = "12"

# Row 49
$ws.Range("D49").Value = "0.002279"
$ws.Range("E49").Value = "-5.09%"
$ws.Range("G49").Value = "12"

# Row 50
$ws.Range("D50").Value = "0.00002100"
$ws.Range("E50").Value = "-0.06%"
$ws.Range("G50").Value = "12"

# Row 51
$ws.Range("D51").Value = "0.0002000"
$ws.Range("E51").Value = "-0.06%"
$ws.Range("G51").Value = "12"
